# Auto-generated script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.641.84'
$ws.Range('E2').Value = '  -7.62%  '
$ws.Range('D3').Value = '2.553.12'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'295.91"
$ws.Range('E5').Value = '  -5.36%  '
$ws.Range('D6').Value = "'91.17"
$ws.Range('E6').Value = '  -7.62%  '
$ws.Range('E7').Value = '  -4.52%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.547"
$ws.Range('E9').Value = '  -5.90%  '
$ws.Range('D10').Value = "'35.49"
$ws.Range('E10').Value = '  -9.06%  '
$ws.Range('E11').Value = '  -4.34%  '
$ws.Range('D12').Value = "'7.61"
$ws.Range('E12').Value = '  -6.46%  '
$ws.Range('D13').Value = '2.939.44'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = '2.539.36'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').Value = "'0.864"
$ws.Range('E16').Value = '  -5.87%  '
$ws.Range('D17').Value = "'14.05"
$ws.Range('E17').Value = '  -5.38%  '
$ws.Range('D18').Value = '42.671.31'
$ws.Range('E18').Value = '  -7.69%  '
$ws.Range('D19').Value = "'6.64"
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('D20').Value = '0.0₃0967'
$ws.Range('E20').Value = '  -5.08%  '
$ws.Range('D21').Value = "'12.48"
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').Value = "'72.45"
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = "'258.91"
$ws.Range('E23').Value = '  -11.60%  '
$ws.Range('E24').Value = '  -6.38%  '
$ws.Range('D25').Value = "'29.62"
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('D26').Value = "'2.12"
$ws.Range('E26').Value = '  -5.88%  '
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = "'9.98"
$ws.Range('E28').Value = '  -7.05%  '
$ws.Range('E29').Value = '  -4.56%  '
$ws.Range('D30').Value = "'35.98"
$ws.Range('E30').Value = '  -5.10%  '
$ws.Range('D31').Value = "'5.89"
$ws.Range('E31').Value = '  -5.45%  '
$ws.Range('D32').Value = "'150.54"
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('E34').Value = '  -5.12%  '
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('D36').Value = "'0.0791"
$ws.Range('E36').Value = '  -5.68%  '
$ws.Range('E37').Value = '  -7.31%  '
$ws.Range('D38').Value = "'24.55"
$ws.Range('E38').Value = '  +9.79%  '
$ws.Range('E39').Value = '  -3.71%  '
$ws.Range('D40').Value = "'15.91"
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').Value = "'3.43"
$ws.Range('E41').Value = '  -4.23%  '
$ws.Range('D42').Value = "'0.0309"
$ws.Range('E42').Value = '  -6.74%  '
$ws.Range('E43').Value = '  -4.07%  '
$ws.Range('D44').Value = '2.067.64'
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = "'84.28"
$ws.Range('E46').Value = '  -13.97%  '
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').Value = '2.797.37'
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = "'1.71"
$ws.Range('E49').Value = '  -2.66%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = "'8.73"
$ws.Range('E50').Value = '  -9.62%  '
$ws.Range('D51').Value = "'103.17"
$ws.Range('E51').Value = '  -5.23%  '
